# Auto-generated edit script: updates market-price-derived columns (H:N)
# across 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the
# scheduled market-data refresh described in the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 5722.1904
$ws.Range("I80").Value = 9272.5
$ws.Range("K80").Value = 27817.5
$ws.Range("M80").Value = -26819.5
$ws.Range("H83").Value = 5722.1904
$ws.Range("I83").Value = 9272.5
$ws.Range("K83").Value = 83452.5
$ws.Range("M83").Value = -78460.5
$ws.Range("H88").Value = 1419
$ws.Range("I88").Value = 402
$ws.Range("J88").Value = 1564.2858
$ws.Range("K88").Value = 402
$ws.Range("L88").Value = 1564.2858
$ws.Range("M88").Value = 4
$ws.Range("N88").Value = -2376.2858
$ws.Range("H91").Value = 1419
$ws.Range("I91").Value = 402
$ws.Range("J91").Value = 1564.2858
$ws.Range("K91").Value = 402
$ws.Range("L91").Value = 1564.2858
$ws.Range("M91").Value = 1002
$ws.Range("N91").Value = -4372.2858
$ws.Range("H94").Value = 2577.0833
$ws.Range("I94").Value = 1780.6666
$ws.Range("J94").Value = 4966.3335
$ws.Range("K94").Value = 1780.6666
$ws.Range("L94").Value = 4966.3335
$ws.Range("M94").Value = -1329.6666
$ws.Range("N94").Value = -5868.3335
$ws.Range("H99").Value = 1417.5454
$ws.Range("J99").Value = 876
$ws.Range("L99").Value = 2628
$ws.Range("N99").Value = -5624
$ws.Range("H137").Value = 2249.5881
$ws.Range("I137").Value = 2370.9
$ws.Range("J137").Value = 2076.2856
$ws.Range("K137").Value = 7112.700000000001
$ws.Range("L137").Value = 6228.8568
$ws.Range("M137").Value = -4562.700000000001
$ws.Range("N137").Value = -11328.8568
$ws.Range("H138").Value = 6479.408
$ws.Range("J138").Value = 7253.946
$ws.Range("L138").Value = 21761.838
$ws.Range("N138").Value = -32041.838
$ws.Range("H141").Value = 4492.9
$ws.Range("I141").Value = 4492.9
$ws.Range("K141").Value = 13478.7
$ws.Range("M141").Value = -8298.699999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5822.016
$ws.Range("I32").Value = 5582.7666
$ws.Range("J32").Value = 12999.5
$ws.Range("K32").Value = 5582.7666
$ws.Range("L32").Value = 12999.5
$ws.Range("M32").Value = -5295.7666
$ws.Range("N32").Value = -13573.5
$ws.Range("H61").Value = 3912.3333
$ws.Range("I61").Value = 3912.3333
$ws.Range("K61").Value = 3912.3333
$ws.Range("M61").Value = -3700.3333
$ws.Range("H74").Value = 1481.625
$ws.Range("I74").Value = 1557.3846
$ws.Range("K74").Value = 1557.3846
$ws.Range("M74").Value = -683.3846000000001
$ws.Range("H77").Value = 1481.625
$ws.Range("I77").Value = 1557.3846
$ws.Range("K77").Value = 7786.923000000001
$ws.Range("M77").Value = -3418.923000000001
$ws.Range("H132").Value = 3352.923
$ws.Range("I132").Value = 3438.3333
$ws.Range("J132").Value = 2883.1667
$ws.Range("K132").Value = 10314.9999
$ws.Range("L132").Value = 8649.500100000001
$ws.Range("M132").Value = -7784.999899999999
$ws.Range("N132").Value = -13709.5001
$ws.Range("H136").Value = 3912.3333
$ws.Range("I136").Value = 3912.3333
$ws.Range("K136").Value = 11736.9999
$ws.Range("M136").Value = -9186.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2685.8276
$ws.Range("I86").Value = 2789.65
$ws.Range("J86").Value = 2455.111
$ws.Range("K86").Value = 2789.65
$ws.Range("L86").Value = 2455.111
$ws.Range("M86").Value = -1666.65
$ws.Range("N86").Value = -4701.111
$ws.Range("H89").Value = 2685.8276
$ws.Range("I89").Value = 2789.65
$ws.Range("J89").Value = 2455.111
$ws.Range("K89").Value = 13948.25
$ws.Range("L89").Value = 12275.555
$ws.Range("M89").Value = -8332.25
$ws.Range("N89").Value = -23507.555
$ws.Range("H94").Value = 1129
$ws.Range("I94").Value = 1274.0476
$ws.Range("J94").Value = 693.8570999999999
$ws.Range("K94").Value = 1274.0476
$ws.Range("L94").Value = 693.8570999999999
$ws.Range("M94").Value = -823.0476000000001
$ws.Range("N94").Value = -1595.8571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N48").ClearContents()
$ws.Range("H31").Value = 48440.547
$ws.Range("I31").Value = 2137.125
$ws.Range("K31").Value = 2137.125
$ws.Range("M31").Value = -1842.125
$ws.Range("H34").Value = 48440.547
$ws.Range("I34").Value = 2137.125
$ws.Range("K34").Value = 2137.125
$ws.Range("M34").Value = -1935.125
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("H105").Value = 642.1111
$ws.Range("I105").Value = 552.9375
$ws.Range("K105").Value = 552.9375
$ws.Range("M105").Value = 1194.0625
$ws.Range("H122").Value = 4823.3335
$ws.Range("I122").Value = 4260.273
$ws.Range("J122").Value = 6371.75
$ws.Range("K122").Value = 12780.819
$ws.Range("L122").Value = 19115.25
$ws.Range("M122").Value = -10330.819
$ws.Range("N122").Value = -24015.25
$ws.Range("H134").Value = 272270.3
$ws.Range("I134").Value = 2055.5293
$ws.Range("K134").Value = 6166.5879
$ws.Range("M134").Value = -3631.5879

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1250
$ws.Range("I12").Value = 2000
$ws.Range("J12").Value = 500
$ws.Range("K12").Value = 6000
$ws.Range("L12").Value = 1500
$ws.Range("M12").Value = -5827
$ws.Range("N12").Value = -1846
$ws.Range("H68").Value = 1819285.4
$ws.Range("J68").Value = 1112362.5
$ws.Range("L68").Value = 3337087.5
$ws.Range("N68").Value = -3338709.5
$ws.Range("H71").Value = 1819285.4
$ws.Range("J71").Value = 1112362.5
$ws.Range("L71").Value = 10011262.5
$ws.Range("N71").Value = -10019374.5
$ws.Range("H140").Value = 5035.9033
$ws.Range("I140").Value = 3248.75
$ws.Range("K140").Value = 9746.25
$ws.Range("M140").Value = -4566.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 11681297
$ws.Range("I14").Value = 14882337
$ws.Range("K14").Value = 14882337
$ws.Range("M14").Value = -14882169
$ws.Range("H80").Value = 1433309.1
$ws.Range("I80").Value = 1005284
$ws.Range("J80").Value = 2503372
$ws.Range("K80").Value = 1005284
$ws.Range("L80").Value = 2503372
$ws.Range("M80").Value = -1004286
$ws.Range("N80").Value = -2505368
$ws.Range("H83").Value = 1433309.1
$ws.Range("I83").Value = 1005284
$ws.Range("J83").Value = 2503372
$ws.Range("K83").Value = 5026420
$ws.Range("L83").Value = 12516860
$ws.Range("M83").Value = -5021428
$ws.Range("N83").Value = -12526844
$ws.Range("H97").Value = 658.56665
$ws.Range("I97").Value = 643.52
$ws.Range("K97").Value = 643.52
$ws.Range("M97").Value = -147.52
$ws.Range("H102").Value = 2295.7144
$ws.Range("I102").Value = 1380.1724
$ws.Range("J102").Value = 4338.077
$ws.Range("K102").Value = 1380.1724
$ws.Range("L102").Value = 4338.077
$ws.Range("M102").Value = 241.8276000000001
$ws.Range("N102").Value = -7582.077
$ws.Range("H113").Value = 597133.4399999999
$ws.Range("I113").Value = 1430504.2
$ws.Range("K113").Value = 1430504.2
$ws.Range("M113").Value = -1428334.2
$ws.Range("H126").Value = 4374.25
$ws.Range("I126").Value = 2998
$ws.Range("J126").Value = 4570.857
$ws.Range("K126").Value = 8994
$ws.Range("L126").Value = 13712.571
$ws.Range("M126").Value = -6524
$ws.Range("N126").Value = -18652.571
$ws.Range("H136").Value = 30765.305
$ws.Range("J136").Value = 30765.305
$ws.Range("L136").Value = 92295.91500000001
$ws.Range("N136").Value = -97395.91500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4753.3335
$ws.Range("I40").Value = 4133.4194
$ws.Range("K40").Value = 4133.4194
$ws.Range("M40").Value = -3997.4194
$ws.Range("H93").Value = 2935.9
$ws.Range("I93").Value = 2571.6667
$ws.Range("J93").Value = 4392.8335
$ws.Range("K93").Value = 2571.6667
$ws.Range("L93").Value = 4392.8335
$ws.Range("M93").Value = -1323.6667
$ws.Range("N93").Value = -6888.8335
$ws.Range("H100").Value = 3158.1667
$ws.Range("I100").Value = 3033.3333
$ws.Range("J100").Value = 3283
$ws.Range("K100").Value = 3033.3333
$ws.Range("L100").Value = 3283
$ws.Range("M100").Value = -2492.3333
$ws.Range("N100").Value = -4365
$ws.Range("H132").Value = 3642
$ws.Range("J132").Value = 4949.25
$ws.Range("L132").Value = 14847.75
$ws.Range("N132").Value = -19907.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 90749.164
$ws.Range("I62").Value = 206599.4
$ws.Range("J62").Value = 7999
$ws.Range("K62").Value = 206599.4
$ws.Range("L62").Value = 7999
$ws.Range("M62").Value = -205975.4
$ws.Range("N62").Value = -9247
$ws.Range("H65").Value = 90749.164
$ws.Range("I65").Value = 206599.4
$ws.Range("J65").Value = 7999
$ws.Range("K65").Value = 1032997
$ws.Range("L65").Value = 39995
$ws.Range("M65").Value = -1029877
$ws.Range("N65").Value = -46235
$ws.Range("H74").Value = 10322.8
$ws.Range("I74").Value = 6999
$ws.Range("K74").Value = 6999
$ws.Range("M74").Value = -6063
$ws.Range("H77").Value = 10322.8
$ws.Range("I77").Value = 6999
$ws.Range("K77").Value = 20997
$ws.Range("M77").Value = -16317
$ws.Range("H113").Value = 260.46667
$ws.Range("I113").Value = 277.46155
$ws.Range("J113").Value = 150
$ws.Range("K113").Value = 832.38465
$ws.Range("L113").Value = 450
$ws.Range("M113").Value = 1337.61535
$ws.Range("N113").Value = -4790
$ws.Range("H132").Value = 16069.25
$ws.Range("J132").Value = 77146.42999999999
$ws.Range("L132").Value = 231439.29
$ws.Range("N132").Value = -236499.29
